$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2127
$ws.Range("I40").Value = 1977.5714
$ws.Range("J40").Value = 2650
$ws.Range("K40").Value = 1977.5714
$ws.Range("L40").Value = 2650
$ws.Range("M40").Value = -1802.5714
$ws.Range("N40").Value = -3000

$ws.Range("H64").Value = 3133.3333
$ws.Range("I64").Value = 3000
$ws.Range("K64").Value = 3000
$ws.Range("M64").Value = -2752

$ws.Range("H67").Value = 3133.3333
$ws.Range("I67").Value = 3000
$ws.Range("K67").Value = 3000
$ws.Range("M67").Value = -2142

$ws.Range("H74").Value = 5000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 5000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 347681.38
$ws.Range("I2").Value = 463399.16
$ws.Range("K2").Value = 463399.16
$ws.Range("M2").Value = -463286.16

$ws.Range("H74").Value = 920.75
$ws.Range("I74").Value = 546.4838999999999
$ws.Range("J74").Value = 3241.2
$ws.Range("K74").Value = 546.4838999999999
$ws.Range("L74").Value = 3241.2
$ws.Range("M74").Value = 327.5161000000001
$ws.Range("N74").Value = -4989.2

$ws.Range("H77").Value = 920.75
$ws.Range("I77").Value = 546.4838999999999
$ws.Range("J77").Value = 3241.2
$ws.Range("K77").Value = 2732.4195
$ws.Range("L77").Value = 16206
$ws.Range("M77").Value = 1635.5805
$ws.Range("N77").Value = -24942

$ws.Range("H102").Value = 1367.3334
$ws.Range("I102").Value = 1052
$ws.Range("J102").Value = 1525
$ws.Range("K102").Value = 1052
$ws.Range("L102").Value = 1525
$ws.Range("M102").Value = 570
$ws.Range("N102").Value = -4769

$ws.Range("H110").Value = 464
$ws.Range("I110").Value = 464
$ws.Range("K110").Value = 464
$ws.Range("M110").Value = 1581

$ws.Range("H116").Value = 347681.38
$ws.Range("I116").Value = 463399.16
$ws.Range("K116").Value = 463399.16
$ws.Range("M116").Value = -461105.16

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 347681.38
$ws.Range("I3").Value = 463399.16
$ws.Range("K3").Value = 463399.16
$ws.Range("M3").Value = -463285.16

$ws.Range("H99").Value = 1542
$ws.Range("I99").Value = 1313.5
$ws.Range("J99").Value = 1999
$ws.Range("K99").Value = 1313.5
$ws.Range("L99").Value = 1999
$ws.Range("M99").Value = 184.5
$ws.Range("N99").Value = -4995

$ws.Range("H134").Value = 4933.759
$ws.Range("I134").Value = 4998.3335
$ws.Range("J134").Value = 4062
$ws.Range("K134").Value = 14995.0005
$ws.Range("L134").Value = 12186
$ws.Range("M134").Value = -12460.0005
$ws.Range("N134").Value = -17256

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10000
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10224

$ws.Range("H58").Value = 2072315.9
$ws.Range("I58").Value = 2072315.9
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2072315.9
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2072112.9
$ws.Range("N58").ClearContents()

$ws.Range("H99").Value = 2637.111
$ws.Range("I99").Value = 2470
$ws.Range("K99").Value = 2470
$ws.Range("M99").Value = -972

$ws.Range("H126").Value = 2637.111
$ws.Range("I126").Value = 2470
$ws.Range("K126").Value = 7410
$ws.Range("M126").Value = -4940

$ws.Range("H132").Value = 1781.7
$ws.Range("I132").Value = 1424.2222
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4272.6666
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1742.6666
$ws.Range("N132").Value = -20057

$ws.Range("H136").Value = 2072315.9
$ws.Range("I136").Value = 2072315.9
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6216947.699999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6214397.699999999
$ws.Range("N136").ClearContents()

$ws.Range("H141").Value = 65500
$ws.Range("J141").Value = 59000
$ws.Range("L141").Value = 59000
$ws.Range("N141").Value = -69360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 269.2
$ws.Range("I6").Value = 269.2
$ws.Range("K6").Value = 807.5999999999999
$ws.Range("M6").Value = -694.5999999999999

$ws.Range("H107").Value = 404.3846
$ws.Range("J107").Value = 368.81818
$ws.Range("L107").Value = 1106.45454
$ws.Range("N107").Value = -4946.45454

$ws.Range("H122").Value = 645.1429000000001
$ws.Range("J122").Value = 646.55554
$ws.Range("L122").Value = 5818.99986
$ws.Range("N122").Value = -10718.99986

$ws.Range("H131").Value = 27652.885
$ws.Range("I131").Value = 765
$ws.Range("J131").Value = 29893.541
$ws.Range("K131").Value = 2295
$ws.Range("L131").Value = 89680.62300000001
$ws.Range("M131").Value = 2745
$ws.Range("N131").Value = -99760.62300000001

$ws.Range("H132").Value = 2230
$ws.Range("J132").Value = 3750
$ws.Range("L132").Value = 33750
$ws.Range("N132").Value = -38810

$ws.Range("H141").Value = 2992.1667
$ws.Range("I141").Value = 2754.4707
$ws.Range("K141").Value = 8263.4121
$ws.Range("M141").Value = -3083.4121

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 20000
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 20000
$ws.Range("N39").Value = -21064

$ws.Range("H126").Value = 3537915
$ws.Range("I126").Value = 4276897
$ws.Range("J126").Value = 335658.34
$ws.Range("K126").Value = 12830691
$ws.Range("L126").Value = 1006975.02
$ws.Range("M126").Value = -12828221
$ws.Range("N126").Value = -1011915.02

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 438333.34
$ws.Range("J2").Value = 130000
$ws.Range("L2").Value = 130000
$ws.Range("N2").Value = -130224

$ws.Range("H29").Value = 10280
$ws.Range("J29").Value = 10280
$ws.Range("L29").Value = 10280
$ws.Range("N29").Value = -10870

$ws.Range("H100").Value = 1111.75
$ws.Range("J100").Value = 1039
$ws.Range("L100").Value = 1039
$ws.Range("N100").Value = -2121

$ws.Range("H132").Value = 2731.6128
$ws.Range("J132").Value = 2940.625
$ws.Range("L132").Value = 8821.875
$ws.Range("N132").Value = -13881.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2400
$ws.Range("J5").Value = 2400
$ws.Range("L5").Value = 2400
$ws.Range("N5").Value = -2624
